{"js": "// HTML ASSIGNMENT 3 COMPLETE\n//\n// 1. Remove the stray \"jm\" run from the very start of the document.\n// 2. Append answers 13-18 after the existing \"12. ...\" paragraph (the\n//    content that used to be the final, empty paragraph becomes the new\n//    last paragraph holding \"...>\" ).\n// 3. Relocate the \"_GoBack\" bookmark so it once again sits at the very end\n//    of the document content, between \"</h1\" and the closing \">\".\n\nconst body = context.document.body;\n\n// --- 1. Drop the leading \"jm\" run -----------------------------------\nconst jm = body.search(\"jm\", { matchCase: true, matchWholeWord: false });\njm.load(\"items\");\nawait context.sync();\nif (jm.items.length > 0) {\n  jm.items[0].delete();\n  await context.sync();\n}\n\n// --- 2. Remove the existing \"_GoBack\" bookmark so we can re-add it in\n//        its new location later without leaving a stale duplicate.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. Locate paragraph 12 and the trailing empty paragraph --------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst trailingEmpty = items[items.length - 1];\nconst para12 = items[items.length - 2];\n\n// --- 4. Insert the new answer paragraphs (13-17) right after #12 ----\nconst p13 = para12.insertParagraph(\n  \"13. You could apply a float:left style to both and then on one of them you could put a style margin-left:30px; or margin-right:30px; so that it creates a barrier and the second div is made 30px away from the first div. You could also do this by having one parent div that surrounds the two child divs and then work with padding through there\",\n  \"After\"\n);\nawait context.sync();\n\nconst p14 = p13.insertParagraph(\n  \"14. the clear:left style makes it so that no content will appear or float beside this element\",\n  \"After\"\n);\nawait context.sync();\n\nconst p15 = p14.insertParagraph(\n  \"15. You would use padding-top:12px; \",\n  \"After\"\n);\nawait context.sync();\n\nconst p16 = p15.insertParagraph(\n  \"16. Relative positioning occurs so that content appears side by side so it will fit relative to the content loaded before and after it. Absolute positioning is just that it is absolute so it has a position that is set on the page.\",\n  \"After\"\n);\nawait context.sync();\n\nconst p17 = p16.insertParagraph(\n  \"17. The z-index style controls which content overlaps where and where the content is visible in reference to the z-axis;\",\n  \"After\"\n);\nawait context.sync();\n\nconst p18a = p17.insertParagraph(\n  \"18. <span style=\\u201dposition:absolute;left:0px;top:0px;\\u201d> Where would you like to</span>\",\n  \"After\"\n);\nawait context.sync();\n\n// --- 5. Turn the old trailing empty paragraph into the final \"#18b\"\n//        paragraph (keeps the document from ending with a stray blank\n//        paragraph, matching the target structure).\ntrailingEmpty.insertText(\n  \"<h1 style=\\u201dposition:absolute;top:20px;left;80px;\\u201d>Go Today?</h1>\",\n  \"Start\"\n);\nawait context.sync();\n\n// --- 6. Re-insert \"_GoBack\" collapsed right before the final \">\" ----\nconst closers = trailingEmpty.search(\">\", { matchCase: true });\nclosers.load(\"items,text\");\nawait context.sync();\n\nconst closer = closers.items[closers.items.length - 1];\nconst beforeCloser = closer.getRange(\"Start\");\nbeforeCloser.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# HTML ASSIGNMENT 3 COMPLETE\n#\n# 1. Remove the stray \"jm\" run from the very start of the document.\n# 2. Append answers 13-18 after the existing \"12. ...\" paragraph (the\n#    content that used to be the final, empty paragraph becomes the new\n#    last paragraph holding \"...>\").\n# 3. Relocate the \"_GoBack\" bookmark so it once again sits at the very end\n#    of the document content, between \"</h1\" and the closing \">\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the leading \"jm\" run ------------------------------------\n$findRng = $d.Content\nif ($findRng.Find.Execute(\"jm\")) {\n  $findRng.Delete()\n}\n\n# --- 2. Remove the existing \"_GoBack\" bookmark so it can be re-added in\n#        its new location later without leaving a stale duplicate.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- 3. Paragraph 12 is (now) the second-to-last paragraph; the very\n#        last paragraph is the trailing empty one.\n$paraCount = $d.Paragraphs.Count\n$para12 = $d.Paragraphs.Item($paraCount - 1)\n\n$answers = @(\n  \"13. You could apply a float:left style to both and then on one of them you could put a style margin-left:30px; or margin-right:30px; so that it creates a barrier and the second div is made 30px away from the first div. You could also do this by having one parent div that surrounds the two child divs and then work with padding through there\",\n  \"14. the clear:left style makes it so that no content will appear or float beside this element\",\n  \"15. You would use padding-top:12px; \",\n  \"16. Relative positioning occurs so that content appears side by side so it will fit relative to the content loaded before and after it. Absolute positioning is just that it is absolute so it has a position that is set on the page.\",\n  \"17. The z-index style controls which content overlaps where and where the content is visible in reference to the z-axis;\",\n  \"18. <span style=\u201dposition:absolute;left:0px;top:0px;\u201d> Where would you like to</span>\"\n)\n\n$prevPara = $para12\nforeach ($answerText in $answers) {\n  $prevPara.Range.InsertParagraphAfter()\n  $paraCount = $d.Paragraphs.Count\n  $newPara = $d.Paragraphs.Item($paraCount - 1)\n  $newPara.Range.Text = $answerText\n  $prevPara = $newPara\n}\n\n# --- 4. Turn the old trailing empty paragraph into the final \"#18b\"\n#        paragraph (keeps the document from ending with a stray blank\n#        paragraph, matching the target structure).\n$lastParaIndex = $d.Paragraphs.Count\n$finalPara = $d.Paragraphs.Item($lastParaIndex)\n$finalPara.Range.Text = \"<h1 style=\u201dposition:absolute;top:20px;left;80px;\u201d>Go Today?</h1>\"\n\n# --- 5. Re-insert \"_GoBack\" collapsed right before the final \">\".\n#        A genuinely zero-length Range is unreliable for Bookmarks.Add in\n#        this host, so insert a throw-away marker character, bookmark the\n#        (non-collapsed) range around it, then delete the marker \u2014 the\n#        bookmark collapses cleanly to that exact position.\n$searchRng = $finalPara.Range.Duplicate()\n$searchRng.Find.Execute(\"</h1>\") | Out-Null\n$splitPos = $searchRng.End - 1\n\n$markerInsertRng = $d.Range($splitPos, $splitPos)\n$markerInsertRng.InsertBefore(\"@\")\n\n$markerRng = $d.Range($splitPos, $splitPos + 1)\n$d.Bookmarks.Add(\"_GoBack\", $markerRng)\n\n$markerRng2 = $d.Range($splitPos, $splitPos + 1)\n$markerRng2.Delete()\n"}
